$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.162.77'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = '2.574.27'
$ws.Range('E3').Value = '  -2.44%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'516.84"
$ws.Range('E5').Value = '  -2.48%  '
$ws.Range('D6').Value = "'138.10"
$ws.Range('E6').Value = '  -5.19%  '
$ws.Range('D7').Value = "'0.997"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('D9').Value = '2.591.50'
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').Value = "'0.0992"
$ws.Range('E11').Value = '  -4.35%  '
$ws.Range('E12').Value = '  -2.48%  '
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '3.038.25'
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').Value = '58.137.62'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').Value = "'20.23"
$ws.Range('E16').Value = '  -2.56%  '
$ws.Range('D17').Value = '2.590.50'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D19').Value = "'335.66"
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('D21').Value = "'10.07"
$ws.Range('E21').Value = '  -4.63%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = "'66.01"
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('E26').Value = '  -3.54%  '
$ws.Range('D27').Value = "'0.996"
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = "'6.96"
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '0.0₃0708'
$ws.Range('E30').Value = '  -10.86%  '
$ws.Range('D31').Value = "'5.86"
$ws.Range('E31').Value = '  -7.65%  '
$ws.Range('D32').Value = "'18.65"
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('D34').Value = "'149.01"
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('D35').Value = "'3.90"
$ws.Range('E35').Value = '  -6.28%  '
$ws.Range('D36').Value = "'1.11"
$ws.Range('E36').Value = '  -5.32%  '
$ws.Range('D37').Value = "'36.17"
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('D38').Value = "'0.827"
$ws.Range('E38').Value = '  -3.32%  '
$ws.Range('D39').Value = "'0.830"
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('D41').Value = "'3.48"
$ws.Range('E41').Value = '  -3.59%  '
$ws.Range('D42').Value = "'0.997"
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = "'271.31"
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('D44').Value = "'10.72"
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = "'0.589"
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('D46').Value = "'0.0942"
$ws.Range('E46').Value = '  -3.79%  '
$ws.Range('D47').Value = "'0.0514"
$ws.Range('E47').Value = '  -3.40%  '
$ws.Range('D48').Value = "'18.38"
$ws.Range('E48').Value = '  -4.53%  '
$ws.Range('D49').Value = '1.969.87'
$ws.Range('E49').Value = '  -3.21%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'4.51"
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = "'0.0217"
$ws.Range('E51').Value = '  -5.08%  '
